$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9598733014176446
$ws.Range("C2").Value = 0.1248264984811058
$ws.Range("E2").Value = 0.1687470901427996
$ws.Range("F2").Value = 2.697307377516523
$ws.Range("G2").Value = 1.617195584468732
$ws.Range("H2").Value = 1.457778898511492
$ws.Range("J2").Value = 0.1139369751323933
$ws.Range("K2").Value = 0.4807606763498313
$ws.Range("L2").Value = 0.3385925687311584
$ws.Range("M2").Value = 0.2846770866607642
$ws.Range("N2").Value = 2.80662303800154
$ws.Range("B3").Value = 0.9260911653678932
$ws.Range("C3").Value = 0.1233853527657658
$ws.Range("E3").Value = 0.1687264768431209
$ws.Range("F3").Value = 2.695877125465529
$ws.Range("G3").Value = 1.619758411478983
$ws.Range("H3").Value = 1.463821603607244
$ws.Range("J3").Value = 0.1136226897306134
$ws.Range("K3").Value = 0.4493613367476428
$ws.Range("L3").Value = 0.3351386875223099
$ws.Range("M3").Value = 0.2779653451877202
$ws.Range("N3").Value = 2.828953703234568
$ws.Range("B4").Value = 0.9058167965146708
$ws.Range("C4").Value = 0.1224825949474422
$ws.Range("E4").Value = 0.1687607500397483
$ws.Range("F4").Value = 2.69614891956013
$ws.Range("G4").Value = 1.622091580609734
$ws.Range("H4").Value = 1.468054073360221
$ws.Range("J4").Value = 0.1134287590252221
$ws.Range("K4").Value = 0.4303013611670679
$ws.Range("L4").Value = 0.3331482675718291
$ws.Range("M4").Value = 0.2739750745281988
$ws.Range("N4").Value = 2.843419383354188
$ws.Range("B5").Value = 0.897672991061512
$ws.Range("C5").Value = 0.122110202683082
$ws.Range("E5").Value = 0.1687865566088682
$ws.Range("F5").Value = 2.696549170415665
$ws.Range("G5").Value = 1.623233378348829
$ws.Range("H5").Value = 1.469910272047727
$ws.Range("J5").Value = 0.1133494983291854
$ws.Range("K5").Value = 0.4225896970491334
$ws.Range("L5").Value = 0.3323700393624307
$ws.Range("M5").Value = 0.2723820196153568
$ws.Range("N5").Value = 2.849504035404436
$ws.Range("B6").Value = 0.8963278694484984
$ws.Range("C6").Value = 0.1220480943957156
$ws.Range("E6").Value = 0.1687915580246564
$ws.Range("F6").Value = 2.696633125916179
$ws.Range("G6").Value = 1.623434509982545
$ws.Range("H6").Value = 1.470226434063775
$ws.Range("J6").Value = 0.1133363233691185
$ws.Range("K6").Value = 0.4213125370912252
$ws.Range("L6").Value = 0.3322428049552144
$ws.Range("M6").Value = 0.2721194918042116
$ws.Range("N6").Value = 2.850525849791101
$ws.Range("B7").Value = 0.9057064872227159
$ws.Range("C7").Value = 0.12247759101718
$ws.Range("E7").Value = 0.1687610500843189
$ws.Range("F7").Value = 2.696153144913794
$ws.Range("G7").Value = 1.622106205881877
$ws.Range("H7").Value = 1.468078574388088
$ws.Range("J7").Value = 0.1134276910173426
$ws.Range("K7").Value = 0.4301971341669173
$ws.Range("L7").Value = 0.3331376388013041
$ws.Range("M7").Value = 0.2739534561914319
$ws.Range("N7").Value = 2.843500674694148
$ws.Range("B8").Value = 0.9481283639477738
$ws.Range("C8").Value = 0.1243332961072525
$ws.Range("E8").Value = 0.168730265537981
$ws.Range("F8").Value = 2.696575735757818
$ws.Range("G8").Value = 1.61792159765244
$ws.Range("H8").Value = 1.45975410343479
$ws.Range("J8").Value = 0.1138288134545355
$ws.Range("K8").Value = 0.4698888125984411
$ws.Range("L8").Value = 0.3373747045789131
$ws.Range("M8").Value = 0.282335813361076
$ws.Range("N8").Value = 2.814166061549869
$ws.Range("B9").Value = 1.035014988773753
$ws.Range("C9").Value = 0.1278311000915053
$ws.Range("E9").Value = 0.1690406591050291
$ws.Range("F9").Value = 2.706517797692683
$ws.Range("G9").Value = 1.615743766124737
$ws.Range("H9").Value = 1.447568956861161
$ws.Range("J9").Value = 0.11460746821637
$ws.Range("K9").Value = 0.5494568362924213
$ws.Range("L9").Value = 0.3467125063339864
$ws.Range("M9").Value = 0.2998065491469788
$ws.Range("N9").Value = 2.762626032621029
$ws.Range("B10").Value = 1.101091067534242
$ws.Range("C10").Value = 0.1303160046281846
$ws.Range("E10").Value = 0.169492737525367
$ws.Range("F10").Value = 2.719367378505567
$ws.Range("G10").Value = 1.617822699711468
$ws.Range("H10").Value = 1.441134636074082
$ws.Range("J10").Value = 0.1151742937060405
$ws.Range("K10").Value = 0.6089680351103084
$ws.Range("L10").Value = 0.354194919072043
$ws.Range("M10").Value = 0.313267495362318
$ws.Range("N10").Value = 2.728404357496011
$ws.Range("B11").Value = 1.13163489839846
$ws.Range("C11").Value = 0.1314282702740073
$ws.Range("E11").Value = 0.1697466354257884
$ws.Range("F11").Value = 2.726415136283819
$ws.Range("G11").Value = 1.619568374007187
$ws.Range("H11").Value = 1.438753264448451
$ws.Range("J11").Value = 0.1154309332781054
$ws.Range("K11").Value = 0.6362693850662424
$ws.Range("L11").Value = 0.3577328469339136
$ws.Range("M11").Value = 0.3195261315836078
$ws.Range("N11").Value = 2.713626657672513
$ws.Range("B12").Value = 1.14327049346133
$ws.Range("C12").Value = 0.1318468657920064
$ws.Range("E12").Value = 0.1698496806148313
$ws.Range("F12").Value = 2.729256604618541
$ws.Range("G12").Value = 1.620344496896905
$ws.Range("H12").Value = 1.437929871392626
$ws.Range("J12").Value = 0.1155279335491386
$ws.Range("K12").Value = 0.64664050304836
$ws.Range("L12").Value = 0.359091746654201
$ws.Range("M12").Value = 0.321915446059009
$ws.Range("N12").Value = 2.708144299356704
$ws.Range("B13").Value = 1.140761485348094
$ws.Range("C13").Value = 0.1317568289998761
$ws.Range("E13").Value = 0.1698271817341208
$ws.Range("F13").Value = 2.728636969870209
$ws.Range("G13").Value = 1.620172226337203
$ws.Range("H13").Value = 1.438103719087493
$ws.Range("J13").Value = 0.1155070511159764
$ws.Range("K13").Value = 0.6444054481267756
$ws.Range("L13").Value = 0.3587982330016644
$ws.Range("M13").Value = 0.3214000079926365
$ws.Range("N13").Value = 2.709319968530139
$ws.Range("B14").Value = 1.132590779631613
$ws.Range("C14").Value = 0.1314627603264995
$ws.Range("E14").Value = 0.1697549749554135
$ws.Range("F14").Value = 2.726645447271963
$ws.Range("G14").Value = 1.619629919613132
$ws.Range("H14").Value = 1.438683952998247
$ws.Range("J14").Value = 0.1154389172801906
$ws.Range("K14").Value = 0.6371219697562367
$ws.Range("L14").Value = 0.3578442610939447
$ws.Range("M14").Value = 0.3197223157973497
$ws.Range("N14").Value = 2.713173343037575
$ws.Range("B15").Value = 1.127594997616143
$ws.Range("C15").Value = 0.1312822971302694
$ws.Range("E15").Value = 0.1697116436187827
$ws.Range("F15").Value = 2.725448056377047
$ws.Range("G15").Value = 1.619312728342223
$ws.Range("H15").Value = 1.43904956812942
$ws.Range("J15").Value = 0.1153971591592402
$ws.Range("K15").Value = 0.6326648771778878
$ws.Range("L15").Value = 0.3572624176339332
$ws.Range("M15").Value = 0.3186971909450662
$ws.Range("N15").Value = 2.715548446403595
$ws.Range("B16").Value = 1.099104682595623
$ws.Range("C16").Value = 0.1302429516329937
$ws.Range("E16").Value = 0.1694771122629213
$ws.Range("F16").Value = 2.718930963520549
$ws.Range("G16").Value = 1.617724716459506
$ws.Range("H16").Value = 1.44130123691663
$ws.Range("J16").Value = 0.1151574965223681
$ws.Range("K16").Value = 0.6071884183635916
$ws.Range("L16").Value = 0.3539663956413079
$ws.Range("M16").Value = 0.3128611875127874
$ws.Range("N16").Value = 2.729386013961971
$ws.Range("B17").Value = 1.081750784406296
$ws.Range("C17").Value = 0.1296007116642954
$ws.Range("E17").Value = 0.1693455645186432
$ws.Range("F17").Value = 2.71524071341581
$ws.Range("G17").Value = 1.6169554397067
$ws.Range("H17").Value = 1.442822254712809
$ws.Range("J17").Value = 0.1150101541981634
$ws.Range("K17").Value = 0.5916179798198584
$ws.Range("L17").Value = 0.351978659826031
$ws.Range("M17").Value = 0.3093155165737542
$ws.Range("N17").Value = 2.738077261551609
$ws.Range("B18").Value = 1.071815007973441
$ws.Range("C18").Value = 0.1292296050475557
$ws.Range("E18").Value = 0.1692744450655503
$ws.Range("F18").Value = 2.713231372658726
$ws.Range("G18").Value = 1.616588265272867
$ws.Range("H18").Value = 1.44374846622587
$ws.Range("J18").Value = 0.1149252932443616
$ws.Range("K18").Value = 0.5826838861561896
$ws.Range("L18").Value = 0.3508480002506076
$ws.Range("M18").Value = 0.307288872647149
$ws.Range("N18").Value = 2.743150584445395
$ws.Range("B19").Value = 1.068458793781105
$ws.Range("C19").Value = 0.1291036611589647
$ws.Range("E19").Value = 0.1692511467224094
$ws.Range("F19").Value = 2.712570493772787
$ws.Range("G19").Value = 1.616476876768786
$ws.Range("H19").Value = 1.444070889539574
$ws.Range("J19").Value = 0.1148965415647041
$ws.Range("K19").Value = 0.579662676338387
$ws.Range("L19").Value = 0.3504673520033634
$ws.Range("M19").Value = 0.3066048767097698
$ws.Range("N19").Value = 2.744881092218243
$ws.Range("B20").Value = 1.083593407592048
$ws.Range("C20").Value = 0.129669255806796
$ws.Range("E20").Value = 0.1693590980591502
$ws.Range("F20").Value = 2.715621833699402
$ws.Range("G20").Value = 1.617029537667847
$ws.Range("H20").Value = 1.44265502455093
$ws.Range("J20").Value = 0.1150258508583377
$ws.Range("K20").Value = 0.5932732451185245
$ws.Range("L20").Value = 0.3521889510927565
$ws.Range("M20").Value = 0.3096916425368477
$ws.Range("N20").Value = 2.737144368417162
$ws.Range("B21").Value = 1.134988836942426
$ws.Range("C21").Value = 0.1315492057676764
$ws.Range("E21").Value = 0.1697759968869867
$ws.Range("F21").Value = 2.727225722613582
$ws.Range("G21").Value = 1.619786084838296
$ws.Range("H21").Value = 1.438511397640568
$ws.Range("J21").Value = 0.1154589348954218
$ws.Range("K21").Value = 0.63926041911418
$ws.Range("L21").Value = 0.3581239466871722
$ws.Range("M21").Value = 0.3202145718087408
$ws.Range("N21").Value = 2.712038429095614
$ws.Range("B22").Value = 1.168982474748532
$ws.Range("C22").Value = 0.132762739555055
$ws.Range("E22").Value = 0.1700886627448561
$ws.Range("F22").Value = 2.735815661063228
$ws.Range("G22").Value = 1.622258406394053
$ws.Range("H22").Value = 1.436260124077037
$ws.Range("J22").Value = 0.1157409061072343
$ws.Range("K22").Value = 0.6695060950553113
$ws.Range("L22").Value = 0.362114448948148
$ws.Range("M22").Value = 0.3272043750291118
$ws.Range("N22").Value = 2.696292627879295
$ws.Range("B23").Value = 1.150802653332732
$ws.Range("C23").Value = 0.1321164335385774
$ws.Range("E23").Value = 0.1699181214737102
$ws.Range("F23").Value = 2.731139077931019
$ws.Range("G23").Value = 1.620877493569964
$ws.Range("H23").Value = 1.437419897852266
$ws.Range("J23").Value = 0.1155905141113038
$ws.Range("K23").Value = 0.6533460858648255
$ws.Range("L23").Value = 0.3599744693673301
$ws.Range("M23").Value = 0.3234635397669123
$ws.Range("N23").Value = 2.704635833780266
$ws.Range("B24").Value = 1.08276022924332
$ws.Range("C24").Value = 0.1296382728384629
$ws.Range("E24").Value = 0.1693529654977226
$ws.Range("F24").Value = 2.715449179596348
$ws.Range("G24").Value = 1.616995804067884
$ws.Range("H24").Value = 1.442730468071844
$ws.Range("J24").Value = 0.1150187548701957
$ws.Range("K24").Value = 0.5925248449614742
$ws.Range("L24").Value = 0.3520938406617802
$ws.Range("M24").Value = 0.3095215591936338
$ws.Range("N24").Value = 2.737565890763484
$ws.Range("B25").Value = 1.011115424867825
$ws.Range("C25").Value = 0.1268998505905827
$ws.Range("E25").Value = 0.1689172033927093
$ws.Range("F25").Value = 2.702853658843679
$ws.Range("G25").Value = 1.615687155681186
$ws.Range("H25").Value = 1.450422784814648
$ws.Range("J25").Value = 0.1143977143598782
$ws.Range("K25").Value = 0.5277464639857214
$ws.Range("L25").Value = 0.3440767817058088
$ws.Range("M25").Value = 0.2949701351980991
$ws.Range("N25").Value = 2.775928647892037
